$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.468.59"
$ws.Range("E2").Value = "  -2.46%  "
$ws.Range("D3").Value = "2.562.50"
$ws.Range("E3").Value = "  -3.70%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -1.48%  "
$ws.Range("D9").Value = "2.574.91"
$ws.Range("E9").Value = "  -3.48%  "
$ws.Range("E10").Value = "  -5.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1000"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.73%  "
$ws.Range("E12").Value = "  -3.21%  "
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("D14").Value = "3.017.38"
$ws.Range("E14").Value = "  -3.54%  "
$ws.Range("D15").Value = "57.433.98"
$ws.Range("E15").Value = "  -2.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.38%  "
$ws.Range("D17").Value = "2.560.17"
$ws.Range("E17").Value = "  -3.54%  "
$ws.Range("E18").Value = "  -2.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "334.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("B27").Value = "Polygon"
$ws.Range("C27").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.401"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.92%  "
$ws.Range("D28").Value = "2.678.72"
$ws.Range("E28").Value = "  -3.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.08%  "
$ws.Range("D30").Value = "0.0₃0746"
$ws.Range("E30").Value = "  -7.16%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.58"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "148.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.47%  "
$ws.Range("E36").Value = "  -3.19%  "
$ws.Range("E37").Value = "  -4.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.836"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -9.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.829"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.46%  "
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "267.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0951"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.587"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.77"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0519"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.94%  "
$ws.Range("D50").Value = "1.967.66"
$ws.Range("E50").Value = "  -4.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.17%  "
